$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in S1: add "dans la source" suffix
$ws.Range("S1").Value = "Valeur en décimales de livres tournois dans la source"

# Update F2: "Import" -> "Imports"
$ws.Range("F2").Value = "Imports"

# Move the active cell / view back to the top-left of the sheet (A1/A2)
$ws.Range("A2").Select()
